$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("§ 275.0-2_P2|llm_response", 36, 12, 24, 0, "documents-2024-11-01-1.json"),
    @("§ 275.0-5_P2|llm_response", 21, 13, 8, 0, "documents-2024-11-01-1.json"),
    @("§ 275.0-7_P2|llm_response", 26, 10, 16, 0, "documents-2024-11-01-1.json")
)

$row = 5
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $row++
}
